$d = $word.ActiveDocument

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*editor of choice.*") {
        $empty = $d.Paragraphs.Item($i + 1)
        $empty.Range.Delete()
        break
    }
}
